$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 520
$ws1.Range("F3").Value = 6206
$ws1.Range("F4").Value = 393
$ws1.Range("F5").Value = 89
$ws1.Range("F6").Value = 123
$ws1.Range("F9").Value = 562

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 520
$ws4.Range("F3").Value = 6206
$ws4.Range("F4").Value = 393
$ws4.Range("F6").Value = 89
$ws4.Range("F7").Value = 123
$ws4.Range("F11").Value = 562
